$wb = $excel.ActiveWorkbook

# --- Sheet: Baseline-limited-ports ---
$ws = $wb.Worksheets.Item("Baseline-limited-ports")
$ws.Cells.Item(2,2).Value = 403.8461538461539
$ws.Cells.Item(3,2).Value = 1005.564443207901
$ws.Cells.Item(4,2).Value = 2033.048913453275
$ws.Cells.Item(5,2).Value = 3265.929493563985
$ws.Cells.Item(6,2).Value = 4490.929259706946
$ws.Cells.Item(7,2).Value = 5721.18277145531
$ws.Cells.Item(8,2).Value = 6952.975212518215
$ws.Cells.Item(9,2).Value = 8361.529460239752
$ws.Cells.Item(10,2).Value = 9798.398227727481
$ws.Cells.Item(11,2).Value = 11239.83113733724
$ws.Cells.Item(12,2).Value = 12675.55655747152
$ws.Cells.Item(13,2).Value = 14113.92646310739
$ws.Cells.Item(14,2).Value = 15555.53703237736
$ws.Cells.Item(15,2).Value = 17092.33012812059
$ws.Cells.Item(16,2).Value = 18735.76715941676
$ws.Cells.Item(17,2).Value = 19997.24720132134
$ws.Cells.Item(18,2).Value = 20811.18751161494
$ws.Cells.Item(19,2).Value = 21625.10452476075
$ws.Cells.Item(20,2).Value = 22433.56579188842
$ws.Cells.Item(21,2).Value = 23242.59035584874
$ws.Cells.Item(22,2).Value = 24056.49507526482
$ws.Cells.Item(23,2).Value = 24870.47017283033
$ws.Cells.Item(24,2).Value = 25675.90383556616
$ws.Cells.Item(25,2).Value = 26487.96803865094
$ws.Cells.Item(26,2).Value = 27301.89555844638
$ws.Cells.Item(27,2).Value = 28114.42466507616
$ws.Cells.Item(28,2).Value = 28918.2418792439
$ws.Cells.Item(29,1).Value = 2057
$ws.Cells.Item(29,2).Value = 29733.34572145313
$ws.Cells.Item(30,1).Value = 2058
$ws.Cells.Item(30,2).Value = 30547.29604162795
$ws.Cells.Item(31,1).Value = 2059
$ws.Cells.Item(31,2).Value = 31356.72997522708
$ws.Cells.Item(32,1).Value = 2060
$ws.Cells.Item(32,2).Value = 32162.51626091804
$ws.Cells.Item(33,1).Value = 2061
$ws.Cells.Item(33,2).Value = 32978.721427244
$ws.Cells.Item(34,1).Value = 2062
$ws.Cells.Item(34,2).Value = 33792.6228746632
$ws.Cells.Item(35,1).Value = 2063
$ws.Cells.Item(35,2).Value = 34606.57870284334

# --- Sheet: Baseline-South-CA ---
$ws = $wb.Worksheets.Item("Baseline-South-CA")
$ws.Cells.Item(2,2).Value = 408.9581304771178
$ws.Cells.Item(3,2).Value = 1018.223551505877
$ws.Cells.Item(4,2).Value = 2045.708021751251
$ws.Cells.Item(5,2).Value = 3884.708584548976
$ws.Cells.Item(6,2).Value = 5718.98763398018
$ws.Cells.Item(7,2).Value = 7553.149796348535
$ws.Cells.Item(8,2).Value = 9409.450157970316
$ws.Cells.Item(9,2).Value = 11634.13373766249
$ws.Cells.Item(10,2).Value = 13881.85941608445
$ws.Cells.Item(11,2).Value = 16201.89501145963
$ws.Cells.Item(12,2).Value = 18651.91618999467
$ws.Cells.Item(13,2).Value = 21108.29106167314
$ws.Cells.Item(14,2).Value = 23558.38741141113
$ws.Cells.Item(15,2).Value = 26008.09314700343
$ws.Cells.Item(16,2).Value = 28457.75669533716
$ws.Cells.Item(17,2).Value = 30420.8863699712
$ws.Cells.Item(18,2).Value = 31234.80765656941
$ws.Cells.Item(19,1).Value = 2047
$ws.Cells.Item(19,2).Value = 32048.7782216854
$ws.Cells.Item(20,1).Value = 2048
$ws.Cells.Item(20,2).Value = 32862.67769209328
$ws.Cells.Item(21,1).Value = 2049
$ws.Cells.Item(21,2).Value = 33678.87009849471
$ws.Cells.Item(22,1).Value = 2050
$ws.Cells.Item(22,2).Value = 34492.79940537025

# --- Sheet: Baseline-Central-CA ---
$ws = $wb.Worksheets.Item("Baseline-Central-CA")
$ws.Cells.Item(2,2).Value = 408.9581304771178
$ws.Cells.Item(3,2).Value = 1018.223551505877
$ws.Cells.Item(4,2).Value = 2045.708021751251
$ws.Cells.Item(5,2).Value = 3281.411881252272
$ws.Cells.Item(6,2).Value = 4514.028502171678
$ws.Cells.Item(7,2).Value = 5746.500988051263
$ws.Cells.Item(8,2).Value = 7001.111673184272
$ws.Cells.Item(9,2).Value = 8622.457111520587
$ws.Cells.Item(10,2).Value = 10870.14146510213
$ws.Cells.Item(11,2).Value = 13117.90525630524
$ws.Cells.Item(12,2).Value = 15365.79889716131
$ws.Cells.Item(13,2).Value = 17619.49245708499
$ws.Cells.Item(14,2).Value = 19867.46126914401
$ws.Cells.Item(15,2).Value = 22115.08578061984
$ws.Cells.Item(16,2).Value = 24434.35790656858
$ws.Cells.Item(17,2).Value = 26512.00400623473
$ws.Cells.Item(18,2).Value = 28129.74250700068
$ws.Cells.Item(19,2).Value = 29747.53028628441
$ws.Cells.Item(20,1).Value = 2048
$ws.Cells.Item(20,2).Value = 31365.20085935633
$ws.Cells.Item(21,1).Value = 2049
$ws.Cells.Item(21,2).Value = 32987.36707041616
$ws.Cells.Item(22,1).Value = 2050
$ws.Cells.Item(22,2).Value = 34492.79940537025

# --- Sheet: Expanded-all-ports ---
$ws = $wb.Worksheets.Item("Expanded-all-ports")
$ws.Cells.Item(2,2).Value = 408.9581304771178
$ws.Cells.Item(3,2).Value = 1018.223551505877
$ws.Cells.Item(4,2).Value = 2045.708021751251
$ws.Cells.Item(5,2).Value = 3884.708584548976
$ws.Cells.Item(6,2).Value = 5718.98763398018
$ws.Cells.Item(7,2).Value = 7553.149796348535
$ws.Cells.Item(8,2).Value = 9409.450157970316
$ws.Cells.Item(9,2).Value = 11634.13373766249
$ws.Cells.Item(10,2).Value = 13881.85941608445
$ws.Cells.Item(11,2).Value = 16201.89501145963
$ws.Cells.Item(12,2).Value = 19741.46952342551
$ws.Cells.Item(13,2).Value = 23825.57011782019
$ws.Cells.Item(14,2).Value = 27899.01258629511
$ws.Cells.Item(15,2).Value = 31972.00520798409
$ws.Cells.Item(16,2).Value = 36034.99308419366
$ws.Cells.Item(17,2).Value = 40001.94910813661
$ws.Cells.Item(18,2).Value = 43271.18792272181
$ws.Cells.Item(19,2).Value = 46324.68813972538
$ws.Cells.Item(20,2).Value = 48774.77288226459
$ws.Cells.Item(21,1).Value = 2049
$ws.Cells.Item(21,2).Value = 51241.61957253295
$ws.Cells.Item(22,1).Value = 2050
$ws.Cells.Item(22,2).Value = 53034.44628541584
$ws.Cells.Item(23,1).Value = 2051
$ws.Cells.Item(23,2).Value = 54373.77672209026
